$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A3:D12) in ascending order by column A (time).
# Row 2 already holds the minimum value, so it is unaffected.
$range = $ws.Range("A3:D12")
$key1 = $ws.Range("A3:A12")

$range.Sort($key1, 1)
